# Update countries & provincias Spain
# This script applies the periodic COVID-19 data refresh to the "Pais" sheet:
#  - Updates the "last updated" timestamp banner in A1
#  - Updates case/death/recovery figures for several countries
#  - Re-labels a few rows whose country ranking changed position
#    (Grecia jumped ahead of Tayikistan/Albania; Mozambique jumped ahead
#    of Somalia/Mayotte) while carrying the displaced countries' previous
#    figures down to the next row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 22 de Agosto de 2020 a las 20:24"

# Estados Unidos
$ws.Range("B4").Value = 5818638
$ws.Range("C4").Value = 21911
$ws.Range("D4").Value = 3133274
$ws.Range("E4").Value = 2505751
$ws.Range("G4").Value = 413
$ws.Range("H4").Value = 179613

# India
$ws.Range("B6").Value = 3043203
$ws.Range("C6").Value = 69835
$ws.Range("D6").Value = 2279797
$ws.Range("E6").Value = 706561
$ws.Range("G6").Value = 917
$ws.Range("H6").Value = 56845

# Turquia
$ws.Range("B21").Value = 257032
$ws.Range("C21").Value = 1309
$ws.Range("D21").Value = 236370
$ws.Range("E21").Value = 14560
$ws.Range("G21").Value = 22
$ws.Range("H21").Value = 6102

# Alemania
$ws.Range("B23").Value = 233355
$ws.Range("C23").Value = 334
$ws.Range("E23").Value = 15075
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 9330

# Israel
$ws.Range("B33").Value = 101856
$ws.Range("C33").Value = 1140
$ws.Range("D33").Value = 78644
$ws.Range("E33").Value = 22393
$ws.Range("G33").Value = 10
$ws.Range("H33").Value = 819

# Marruecos
$ws.Range("B53").Value = 50812
$ws.Range("C53").Value = 1565
$ws.Range("D53").Value = 35040
$ws.Range("E53").Value = 14914
$ws.Range("G53").Value = 41
$ws.Range("H53").Value = 858

# Argelia
$ws.Range("B58").Value = 41068
$ws.Range("C58").Value = 401
$ws.Range("D58").Value = 28874
$ws.Range("E58").Value = 10770
$ws.Range("G58").Value = 6
$ws.Range("H58").Value = 1424

# Zambia
$ws.Range("B88").Value = 10831
$ws.Range("C88").Value = 204
$ws.Range("D88").Value = 9942
$ws.Range("E88").Value = 610
$ws.Range("G88").Value = 2
$ws.Range("H88").Value = 279

# Grecia moves ahead of Tayikistan and Albania (rows 96-98 re-ranked)
$ws.Range("A96").Value = "Grecia"
$ws.Range("B96").Value = 8381
$ws.Range("C96").Value = 243
$ws.Range("D96").Value = 3804
$ws.Range("E96").Value = 4337
$ws.Range("G96").Value = 2
$ws.Range("H96").Value = 240

$ws.Range("A97").Value = "Tayikistan"
$ws.Range("B97").Value = 8277
$ws.Range("C97").Value = 36
$ws.Range("D97").Value = 7072
$ws.Range("E97").Value = 1139
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 66

$ws.Range("A98").Value = "Albania"
$ws.Range("B98").Value = 8275
$ws.Range("C98").Value = 156
$ws.Range("D98").Value = 4184
$ws.Range("E98").Value = 3846
$ws.Range("G98").Value = 5
$ws.Range("H98").Value = 245

# Zimbabue
$ws.Range("B105").Value = 5893
$ws.Range("C105").Value = 78
$ws.Range("D105").Value = 4629
$ws.Range("E105").Value = 1111
$ws.Range("G105").Value = 1
$ws.Range("H105").Value = 153

# Mozambique moves ahead of Somalia and Mayotte (rows 122-124 re-ranked)
$ws.Range("A122").Value = "Mozambique"
$ws.Range("B122").Value = 3304
$ws.Range("C122").Value = 109
$ws.Range("D122").Value = 1474
$ws.Range("E122").Value = 1810
$ws.Range("H122").Value = 20

$ws.Range("A123").Value = "Somalia"
$ws.Range("B123").Value = 3265
$ws.Range("D123").Value = 2396
$ws.Range("E123").Value = 776
$ws.Range("H123").Value = 93

$ws.Range("A124").Value = "Mayotte"
$ws.Range("B124").Value = 3237
$ws.Range("D124").Value = 2964
$ws.Range("E124").Value = 234
$ws.Range("H124").Value = 39

# Yemen
$ws.Range("B141").Value = 1907
$ws.Range("C141").Value = 1
$ws.Range("D141").Value = 1066
$ws.Range("E141").Value = 295
$ws.Range("G141").Value = 4
$ws.Range("H141").Value = 546
